$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.960.47"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.642.96"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'215.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'0.5057"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").Value = "'0.2575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "'0.06419"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'19.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "'0.07757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.655.86"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "'4.271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "1.870.90"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "'0.5458"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "0.0₅7938"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "'64.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "25.987.69"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "'201.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").Value = "'4.387"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "'9.908"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "'5.992"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'1.869"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "'141.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").Value = "'0.1138"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "'6.832"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").Value = "'15.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "'1.243"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "'0.04928"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("D32").Value = "'3.273"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "'3.213"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'1.542"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'2.369"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").Value = "'0.8939"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").Value = "'2.621"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").Value = "1.155.77"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "'0.5590"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").Value = "'0.01568"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'1.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("D42").Value = "'5.716"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "'0.8099"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").Value = "'99.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "1.782.40"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "0.0₈118"
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("D47").Value = "'0.4513"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "'54.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "'0.05048"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
